$d = $word.ActiveDocument

# Helper: replace a Range's contents with a literal <w:p>...</w:p> fragment by
# wrapping it as a single-part WordOpenXML package and calling Range.InsertXML --
# this is the supported way to push hand-built run/proofErr markup through the COM
# object model (Range.WordOpenXML itself is read-only).
function Set-RangeOpenXml($Range, $InnerXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $InnerXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($pkg) | Out-Null
}

# --- Paragraph 2 ("${block__customer}" / "Customer: ..." / "Address: ..." / "${/block__customer}") ---
# Re-split the runs and wrap the dictionary-word fragments in w:proofErr spellStart/spellEnd
# markers (cosmetic spell-check bookkeeping only -- the concatenated text is unchanged).
$p2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>block</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>_</w:t></w:r><w:r><w:t>_</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>customer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r><w:r><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Customer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: ${</w:t></w:r><w:r><w:t>block__</w:t></w:r><w:r><w:t>customer</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>name}</w:t></w:r><w:r><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Address</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: ${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>block</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>__</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>customer</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>address</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r><w:r><w:br/><w:t>${/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>block</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>_</w:t></w:r><w:r><w:t>_</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>customer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p>'
Set-RangeOpenXml $d.Paragraphs(2).Range $p2Xml

# --- Table cells: "row__" merge-field prefix renamed to "table__" ---
$tbl = $d.Tables(1)

# Cell 1: ${row__account.id} -> ${table__account.id}
$cell1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>${</w:t></w:r><w:r><w:t>table</w:t></w:r><w:r><w:t>__account.id}</w:t></w:r></w:p>'
Set-RangeOpenXml $tbl.Cell(1,1).Range $cell1Xml

# Cell 2: ${row__account.name} -> ${table__account.name}
$cell2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>${</w:t></w:r><w:r><w:t>table</w:t></w:r><w:r><w:t>__account.name}</w:t></w:r></w:p>'
Set-RangeOpenXml $tbl.Cell(1,2).Range $cell2Xml

# Cell 3: ${row__account.number} -> ${table__account.number} (with proofErr markers)
$cell3Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>__</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>account.number</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p>'
Set-RangeOpenXml $tbl.Cell(1,3).Range $cell3Xml
